# Updated datasets from Daan allowing 500ms gaps
# Recomputes gaze-shift / fixation-duration summary columns (F:K) in the
# LLL gazepath trial table after re-deriving the data with a 500ms gap
# allowance. For a handful of trials this also changes whether a
# first_shift/latency pair (columns F/G) is present at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 1576
$ws.Range("K2").Value = 0
# Row 3
$ws.Range("J3").Value = 503
# Row 4
$ws.Range("J4").Value = 917
# Row 6
$ws.Range("F6").Value = "congruent"
$ws.Range("G6").Value = 8269
$ws.Range("H6").Value = 1
$ws.Range("J6").Value = 1706
$ws.Range("K6").Value = 235
# Row 10
$ws.Range("K10").Value = 0
# Row 12
$ws.Range("J12").Value = 1882
# Row 13
$ws.Range("J13").Value = 1744
$ws.Range("K13").Value = 0
# Row 16
$ws.Range("I16").Value = 0
# Row 21
$ws.Range("H21").Value = 3
# Row 23
$ws.Range("K23").Value = 0
# Row 24
$ws.Range("K24").Value = 0
# Row 25
$ws.Range("J25").Value = 0
# Row 28
$ws.Range("F28").Value = "incongruent"
$ws.Range("G28").Value = 6857
$ws.Range("I28").Value = 1
# Row 31
$ws.Range("J31").Value = 0
# Row 37
$ws.Range("J37").Value = 0
# Row 40
$ws.Range("H40").Value = 4
# Row 43
$ws.Range("J43").Value = 681
# Row 53
$ws.Range("K53").Value = 0
# Row 57
$ws.Range("K57").Value = 0
# Row 59
$ws.Range("F59").Value = $null
$ws.Range("G59").Value = $null
$ws.Range("I59").Value = 0
# Row 68
$ws.Range("J68").Value = 0
# Row 71
$ws.Range("H71").Value = 3
# Row 73
$ws.Range("F73").Value = $null
$ws.Range("G73").Value = $null
$ws.Range("I73").Value = 0
# Row 74
$ws.Range("J74").Value = 227
# Row 77
$ws.Range("H77").Value = 2
# Row 78
$ws.Range("K78").Value = 0
# Row 81
$ws.Range("J81").Value = 0
# Row 83
$ws.Range("J83").Value = 769
# Row 84
$ws.Range("G84").Value = 5297
# Row 85
$ws.Range("J85").Value = 0
# Row 86
$ws.Range("J86").Value = 0
# Row 87
$ws.Range("K87").Value = 0
# Row 93
$ws.Range("F93").Value = "congruent"
$ws.Range("G93").Value = 5423
$ws.Range("H93").Value = 1
$ws.Range("J93").Value = 0
# Row 94
$ws.Range("F94").Value = "congruent"
$ws.Range("G94").Value = 8717
$ws.Range("I94").Value = 0
# Row 95
$ws.Range("J95").Value = 695
# Row 96
$ws.Range("J96").Value = 0
# Row 97
$ws.Range("K97").Value = 0
# Row 98
$ws.Range("J98").Value = 447
# Row 99
$ws.Range("J99").Value = 495
# Row 106
$ws.Range("J106").Value = 2005
# Row 107
$ws.Range("K107").Value = 0
# Row 111
$ws.Range("G111").Value = 6953
$ws.Range("H111").Value = 1
# Row 113
$ws.Range("H113").Value = 0
# Row 116
$ws.Range("J116").Value = 0
# Row 117
$ws.Range("J117").Value = 0
# Row 119
$ws.Range("K119").Value = 0
# Row 122
$ws.Range("K122").Value = 0
# Row 123
$ws.Range("H123").Value = 2
# Row 141
$ws.Range("F141").Value = "incongruent"
$ws.Range("G141").Value = 8211
$ws.Range("I141").Value = 1
# Row 142
$ws.Range("J142").Value = 0
# Row 146
$ws.Range("J146").Value = 0
# Row 150
$ws.Range("J150").Value = 0
# Row 153
$ws.Range("G153").Value = 9725
$ws.Range("I153").Value = 1
# Row 156
$ws.Range("J156").Value = 0
# Row 157
$ws.Range("J157").Value = 423
# Row 158
$ws.Range("J158").Value = 777
$ws.Range("K158").Value = 627
# Row 159
$ws.Range("K159").Value = 0
# Row 161
$ws.Range("K161").Value = 0
# Row 165
$ws.Range("J165").Value = 613
$ws.Range("K165").Value = 0
# Row 167
$ws.Range("H167").Value = 2
$ws.Range("J167").Value = 0
# Row 168
$ws.Range("J168").Value = 211
# Row 170
$ws.Range("K170").Value = 0
# Row 171
$ws.Range("K171").Value = 867
# Row 177
$ws.Range("K177").Value = 341
# Row 178
$ws.Range("K178").Value = 387
# Row 184
$ws.Range("J184").Value = 683
# Row 185
$ws.Range("J185").Value = 0
# Row 186
$ws.Range("F186").Value = "incongruent"
$ws.Range("G186").Value = 6697
$ws.Range("I186").Value = 1
$ws.Range("K186").Value = 601
# Row 187
$ws.Range("J187").Value = 0
$ws.Range("K187").Value = 1285
# Row 189
$ws.Range("K189").Value = 0
# Row 195
$ws.Range("J195").Value = 591
# Row 196
$ws.Range("J196").Value = 1610
# Row 201
$ws.Range("K201").Value = 0
# Row 203
$ws.Range("I203").Value = 1
$ws.Range("J203").Value = 159
# Row 207
$ws.Range("J207").Value = 0
# Row 208
$ws.Range("F208").Value = $null
$ws.Range("G208").Value = $null
$ws.Range("H208").Value = 0
# Row 209
$ws.Range("F209").Value = $null
$ws.Range("G209").Value = $null
$ws.Range("I209").Value = 0
# Row 210
$ws.Range("F210").Value = "congruent"
$ws.Range("G210").Value = 5813
$ws.Range("H210").Value = 1
# Row 212
$ws.Range("K212").Value = 0
# Row 214
$ws.Range("J214").Value = 1264
# Row 216
$ws.Range("J216").Value = 608
# Row 217
$ws.Range("K217").Value = 1235
# Row 218
$ws.Range("K218").Value = 0
# Row 219
$ws.Range("J219").Value = 0
# Row 221
$ws.Range("G221").Value = 9141
$ws.Range("H221").Value = 1
# Row 227
$ws.Range("K227").Value = 0
# Row 229
$ws.Range("J229").Value = 483
# Row 230
$ws.Range("J230").Value = 237
# Row 231
$ws.Range("J231").Value = 0
$ws.Range("K231").Value = 0
# Row 232
$ws.Range("K232").Value = 717
# Row 236
$ws.Range("K236").Value = 0
# Row 237
$ws.Range("J237").Value = 447
# Row 239
$ws.Range("J239").Value = 401
$ws.Range("K239").Value = 0
# Row 240
$ws.Range("J240").Value = 0
# Row 242
$ws.Range("F242").Value = "congruent"
$ws.Range("G242").Value = 4931
$ws.Range("H242").Value = 1
$ws.Range("K242").Value = 0
# Row 247
$ws.Range("J247").Value = 0
# Row 248
$ws.Range("J248").Value = 0
# Row 250
$ws.Range("J250").Value = 511
# Row 251
$ws.Range("J251").Value = 1105
# Row 252
$ws.Range("J252").Value = 411
# Row 253
$ws.Range("J253").Value = 0
$ws.Range("K253").Value = 261
# Row 254
$ws.Range("J254").Value = 0
# Row 255
$ws.Range("J255").Value = 1159
$ws.Range("K255").Value = 1908
# Row 266
$ws.Range("F266").Value = "congruent"
$ws.Range("G266").Value = 5757
$ws.Range("H266").Value = 1
# Row 267
$ws.Range("J267").Value = 0
# Row 268
$ws.Range("J268").Value = 0
# Row 269
$ws.Range("J269").Value = 0
# Row 270
$ws.Range("J270").Value = 0
$ws.Range("K270").Value = 0
# Row 272
$ws.Range("J272").Value = 163
$ws.Range("K272").Value = 0
# Row 273
$ws.Range("K273").Value = 0
# Row 274
$ws.Range("J274").Value = 779
$ws.Range("K274").Value = 1445
# Row 276
$ws.Range("K276").Value = 1256
# Row 277
$ws.Range("J277").Value = 1245
# Row 288
$ws.Range("J288").Value = 0
# Row 290
$ws.Range("J290").Value = 1361
# Row 298
$ws.Range("K298").Value = 0
# Row 299
$ws.Range("J299").Value = 1242
# Row 300
$ws.Range("K300").Value = 965
# Row 310
$ws.Range("J310").Value = 913
# Row 319
$ws.Range("K319").Value = 0
